$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Remove the SUBTOTAL formulas from row 28 (bugfix - E28 becomes blank, G28 becomes a plain value)
$ws.Range("E28").ClearContents()
$ws.Range("G28").Value = 7

# Apply center alignment styling to the relevant cells (new duplicate style, same visual: center aligned)
$ws.Range("E25:F25").HorizontalAlignment = -4108   # xlCenter
$ws.Range("G26:H26").HorizontalAlignment = -4108
$ws.Range("E27:E28").HorizontalAlignment = -4108
$ws.Range("G27").HorizontalAlignment = -4108

# New merges
$ws.Range("G26:H26").Merge()
$ws.Range("E27:E28").Merge()

# Add new data rows
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 4

$ws.Range("E31").Value = 3
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 2

$ws.Range("E32").Value = "SUBTOTAL"
$ws.Range("F32").Value = 3
$ws.Range("G32").Value = 4
$ws.Range("H32").Value = 5

$ws.Range("E33").Value = "TOTAL"
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 7
$ws.Range("H33").Value = 8

# Update view: selection and top-left cell
$ws.Range("J27").Select()
$excel.ActiveWindow.ScrollRow = 18
